# Fitur Product & Kategori Admin
# Adds a "status" column to the Category sheet (between "slug" and
# "created_at"), defaulting every existing category row to status = 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (created_at), pushing created_at /
# updated_at and everything else one column to the right.
$ws.Columns("D:D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "status"

# Default every existing category row to an active status (1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = 1
}

# Match the author's final selection in the sheet.
$ws.Range("E4").Select()
